# Auto-generated cell updates reflecting the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''328.68'
$ws.Range("E2").Value = '''0.24%'
$ws.Range("G2").Value = '''18'
$ws.Range("D3").Value = '''44.32'
$ws.Range("E3").Value = '''1.03%'
$ws.Range("G3").Value = '''18'
$ws.Range("D4").Value = '''5.581'
$ws.Range("E4").Value = '''2.83%'
$ws.Range("G4").Value = '''18'
$ws.Range("D5").Value = '''0.08094'
$ws.Range("E5").Value = '''-0.10%'
$ws.Range("G5").Value = '''18'
$ws.Range("D6").Value = '''1.943'
$ws.Range("E6").Value = '''2.50%'
$ws.Range("G6").Value = '''18'
$ws.Range("D7").Value = '''0.9533'
$ws.Range("E7").Value = '''0.74%'
$ws.Range("G7").Value = '''18'
$ws.Range("D8").Value = '''2.564'
$ws.Range("E8").Value = '''-7.67%'
$ws.Range("G8").Value = '''18'
$ws.Range("D9").Value = '''0.1187'
$ws.Range("E9").Value = '''0.78%'
$ws.Range("G9").Value = '''18'
$ws.Range("D10").Value = '''0.1851'
$ws.Range("E10").Value = '''-1.97%'
$ws.Range("G10").Value = '''18'
$ws.Range("D11").Value = '''0.09808'
$ws.Range("E11").Value = '''1.90%'
$ws.Range("G11").Value = '''18'
$ws.Range("D12").Value = '''0.04493'
$ws.Range("E12").Value = '''6.25%'
$ws.Range("G12").Value = '''18'
$ws.Range("D13").Value = '''0.1068'
$ws.Range("E13").Value = '''-0.13%'
$ws.Range("G13").Value = '''18'
$ws.Range("D14").Value = '''0.001285'
$ws.Range("E14").Value = '''0.90%'
$ws.Range("G14").Value = '''18'
$ws.Range("E15").Value = '''-4.22%'
$ws.Range("G15").Value = '''18'
$ws.Range("D16").Value = '''0.005868'
$ws.Range("E16").Value = '''-3.86%'
$ws.Range("G16").Value = '''18'
$ws.Range("B17").Value = 'HotbitToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
$ws.Range("D17").Value = '''0.004373'
$ws.Range("E17").Value = '''1.59%'
$ws.Range("G17").Value = '''18'
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18").Value = '''3.385'
$ws.Range("E18").Value = '''-5.01%'
$ws.Range("G18").Value = '''18'
$ws.Range("B19").Value = 'GateToken'
$ws.Range("C19").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D19").Value = '''4.312'
$ws.Range("E19").Value = '''-0.19%'
$ws.Range("G19").Value = '''18'
$ws.Range("B20").Value = 'BitpandaEcosystemToken'
$ws.Range("C20").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("D20").Value = '''0.3480'
$ws.Range("E20").Value = '''-1.06%'
$ws.Range("G20").Value = '''18'
$ws.Range("B21").Value = 'MCDex'
$ws.Range("C21").Value = 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
$ws.Range("D21").Value = '''10.16'
$ws.Range("E21").Value = '''16.42%'
$ws.Range("G21").Value = '''18'
$ws.Range("B22").Value = 'ProBitToken'
$ws.Range("C22").Value = 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
$ws.Range("D22").Value = '''0.1420'
$ws.Range("E22").Value = '''4.35%'
$ws.Range("G22").Value = '''18'
$ws.Range("B23").Value = 'ZBToken'
$ws.Range("C23").Value = 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'
$ws.Range("D23").Value = '''0.2505'
$ws.Range("E23").Value = '''-3.91%'
$ws.Range("G23").Value = '''18'
$ws.Range("B24").Value = 'BitKan'
$ws.Range("C24").Value = 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
$ws.Range("D24").Value = '''0.001246'
$ws.Range("E24").Value = '''0.39%'
$ws.Range("G24").Value = '''18'
$ws.Range("E25").Value = '''-4.04%'
$ws.Range("G25").Value = '''18'
$ws.Range("E26").Value = '''-1.00%'
$ws.Range("G26").Value = '''18'
$ws.Range("G27").Value = '''18'
$ws.Range("G28").Value = '''18'
$ws.Range("G29").Value = '''18'
$ws.Range("G30").Value = '''18'
$ws.Range("G31").Value = '''18'
$ws.Range("G32").Value = '''18'
$ws.Range("G33").Value = '''18'
$ws.Range("G34").Value = '''18'
$ws.Range("G35").Value = '''18'
$ws.Range("G36").Value = '''18'
$ws.Range("G37").Value = '''18'
$ws.Range("D38").Value = '''0.02674'
$ws.Range("E38").Value = '''-1.38%'
$ws.Range("G38").Value = '''18'
$ws.Range("D39").Value = '''0.05564'
$ws.Range("E39").Value = '''0.49%'
$ws.Range("G39").Value = '''18'
$ws.Range("D40").Value = '''0.007565'
$ws.Range("E40").Value = '''-2.95%'
$ws.Range("G40").Value = '''18'
$ws.Range("D41").Value = '''0.1409'
$ws.Range("E41").Value = '''0.63%'
$ws.Range("G41").Value = '''18'
$ws.Range("D42").Value = '''0.007966'
$ws.Range("E42").Value = '''-18.67%'
$ws.Range("G42").Value = '''18'
$ws.Range("D43").Value = '''0.002016'
$ws.Range("E43").Value = '''-5.65%'
$ws.Range("G43").Value = '''18'
$ws.Range("D44").Value = '''0.008402'
$ws.Range("E44").Value = '''-12.80%'
$ws.Range("G44").Value = '''18'
$ws.Range("D45").Value = '''0.00007174'
$ws.Range("E45").Value = '''0.89%'
$ws.Range("G45").Value = '''18'
$ws.Range("D46").Value = '''0.00000000750'
$ws.Range("E46").Value = '''-0.83%'
$ws.Range("G46").Value = '''18'
$ws.Range("B47").Value = 'CoinbaseStockToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin'
$ws.Range("D47").Value = '''0.002270'
$ws.Range("E47").Value = '''-0.76%'
$ws.Range("G47").Value = '''18'
$ws.Range("B48").Value = 'BOLO'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ogrGe0dEab+bolo-bolo'
$ws.Range("D48").Value = '''0.003871'
$ws.Range("E48").Value = '''11.29%'
$ws.Range("G48").Value = '''18'
$ws.Range("D49").Value = '''0.00002100'
$ws.Range("E49").Value = '''-0.83%'
$ws.Range("G49").Value = '''18'
$ws.Range("D50").Value = '''0.0002000'
$ws.Range("E50").Value = '''-0.83%'
$ws.Range("G50").Value = '''18'
$ws.Range("G51").Value = '''18'
